# Applies the latest "cryptos list" data refresh to Sheet1.
# Column A (rank index) is untouched; B/C/D/E (Coin/Link/Price/Volume) get
# refreshed numbers, and two coin pairs swapped rank position this cycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume text in this sheet is stored as literal text (e.g. prices use
# "25.837.39" style grouping, not real numbers), never as numeric cells.
# Excel's normal Value coercion would turn number-looking text (e.g.
# "215.05") into a real number (and mangle it via float rounding), so force
# each written cell to Text format first, then flip the format back to
# Normal so the saved style index matches the untouched cells around it.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Subscript digits used by very-small-cap coin prices (e.g. 0.0[sub3]0765)
# are built from code points and stashed in a variable first.
$sub3 = [char]0x2083
$sub6 = [char]0x2086
$d16 = [string]::Concat("0.0", $sub3, "0764")
$d45 = [string]::Concat("0.0", $sub6, "0108")

# --- Column D (Price) updates -----------------------------------------------
Set-TextValue $ws.Range("D2")  "25.837.39"
Set-TextValue $ws.Range("D3")  "1.633.63"
Set-TextValue $ws.Range("D5")  "215.05"
Set-TextValue $ws.Range("D6")  "0.508"
Set-TextValue $ws.Range("D8")  "0.258"
Set-TextValue $ws.Range("D10") "19.89"
Set-TextValue $ws.Range("D11") "0.0782"
Set-TextValue $ws.Range("D13") "1.859.22"
Set-TextValue $ws.Range("D14") "1.632.42"
Set-TextValue $ws.Range("D16") $d16
Set-TextValue $ws.Range("D17") "63.07"
Set-TextValue $ws.Range("D18") "25.838.01"
Set-TextValue $ws.Range("D20") "193.51"
Set-TextValue $ws.Range("D22") "9.91"
Set-TextValue $ws.Range("D26") "138.72"
Set-TextValue $ws.Range("D29") "15.55"
Set-TextValue $ws.Range("D31") "0.0493"
Set-TextValue $ws.Range("D43") "99.36"
Set-TextValue $ws.Range("D44") "0.798"
Set-TextValue $ws.Range("D45") $d45
Set-TextValue $ws.Range("D46") "55.37"
Set-TextValue $ws.Range("D47") "0.422"

# --- Column E (Volume 1h) updates (always plain text already) --------------
$ws.Range("E2").Value  = "  -0.03%  "
$ws.Range("E3").Value  = "  +0.20%  "
$ws.Range("E4").Value  = "  -0.11%  "
$ws.Range("E5").Value  = "  -0.26%  "
$ws.Range("E6").Value  = "  -0.50%  "
$ws.Range("E7").Value  = "  -0.07%  "
$ws.Range("E8").Value  = "  -0.38%  "
$ws.Range("E9").Value  = "  +0.27%  "
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -5.63%  "
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  -4.90%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E51").Value = "  -0.62%  "

# --- Row re-ordering: two coin pairs swapped rank position this refresh ----
# Rows 38/39 swap (Maker <-> ImmutableX) and rows 49/50 swap
# (EnergySwap <-> SynthetixNetwork). The rank column (A) stays as-is; only
# B (Coin), C (Link), D (Price) and E (Volume) move with the coin.

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D38") "0.548"
$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D39") "1.120.88"
$ws.Range("E39").Value = "  -0.97%  "

$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Range("D49") "2.36"
$ws.Range("E49").Value = "  +9.10%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.60"
$ws.Range("E50").Value = "  -0.33%  "
